{"js": "const body = context.document.body;\n\nconst replacements = [\n  [\"2025-04-23 Wednesday\", \"2025-04-24 Thursday\"],\n  [\"792\u00d73=2376\", \"774\u00d77=5418\"],\n  [\"743\u00d76=4458\", \"128\u00d76=768\"],\n  [\"368\u00d75=1840\", \"578\u00d79=5202\"],\n  [\"537\u00d76=3222\", \"367\u00d78=2936\"],\n  [\"110\u00d78=880\", \"458\u00d73=1374\"],\n  [\"212\u00d76=1272\", \"555\u00d72=1110\"],\n  [\"377\u00d78=3016\", \"602\u00d76=3612\"],\n  [\"908\u00d76=5448\", \"493\u00d77=3451\"],\n  [\"851\u00d72=1702\", \"534\u00d75=2670\"],\n  [\"448\u00d74=1792\", \"689\u00d74=2756\"],\n  [\"188\u00d73=564\", \"467\u00d73=1401\"],\n  [\"466\u00d77=3262\", \"101\u00d72=202\"],\n  [\"266\u00d77=1862\", \"280\u00d73=840\"],\n  [\"576\u00d78=4608\", \"764\u00d73=2292\"],\n  [\"446\u00d78=3568\", \"825\u00d75=4125\"],\n  [\"743\u00d74=2972\", \"812\u00d75=4060\"],\n  [\"860\u00d74=3440\", \"208\u00d74=832\"],\n  [\"195\u00d72=390\", \"889\u00d72=1778\"],\n  [\"126\u00d78=1008\", \"475\u00d79=4275\"],\n  [\"586\u00d73=1758\", \"976\u00d78=7808\"],\n  [\"441\u00d73=1323\", \"277\u00d72=554\"],\n  [\"613\u00d76=3678\", \"663\u00d77=4641\"],\n  [\"832\u00d74=3328\", \"910\u00d78=7280\"],\n  [\"824\u00d76=4944\", \"702\u00d77=4914\"],\n  [\"515\u00d74=2060\", \"579\u00d72=1158\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($doc, $old, $new) {\n    $range = $doc.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n\nReplace-Text $d '2025-04-23 Wednesday' '2025-04-24 Thursday'\nReplace-Text $d '792\u00d73=2376' '774\u00d77=5418'\nReplace-Text $d '743\u00d76=4458' '128\u00d76=768'\nReplace-Text $d '368\u00d75=1840' '578\u00d79=5202'\nReplace-Text $d '537\u00d76=3222' '367\u00d78=2936'\nReplace-Text $d '110\u00d78=880' '458\u00d73=1374'\nReplace-Text $d '212\u00d76=1272' '555\u00d72=1110'\nReplace-Text $d '377\u00d78=3016' '602\u00d76=3612'\nReplace-Text $d '908\u00d76=5448' '493\u00d77=3451'\nReplace-Text $d '851\u00d72=1702' '534\u00d75=2670'\nReplace-Text $d '448\u00d74=1792' '689\u00d74=2756'\nReplace-Text $d '188\u00d73=564' '467\u00d73=1401'\nReplace-Text $d '466\u00d77=3262' '101\u00d72=202'\nReplace-Text $d '266\u00d77=1862' '280\u00d73=840'\nReplace-Text $d '576\u00d78=4608' '764\u00d73=2292'\nReplace-Text $d '446\u00d78=3568' '825\u00d75=4125'\nReplace-Text $d '743\u00d74=2972' '812\u00d75=4060'\nReplace-Text $d '860\u00d74=3440' '208\u00d74=832'\nReplace-Text $d '195\u00d72=390' '889\u00d72=1778'\nReplace-Text $d '126\u00d78=1008' '475\u00d79=4275'\nReplace-Text $d '586\u00d73=1758' '976\u00d78=7808'\nReplace-Text $d '441\u00d73=1323' '277\u00d72=554'\nReplace-Text $d '613\u00d76=3678' '663\u00d77=4641'\nReplace-Text $d '832\u00d74=3328' '910\u00d78=7280'\nReplace-Text $d '824\u00d76=4944' '702\u00d77=4914'\nReplace-Text $d '515\u00d74=2060' '579\u00d72=1158'\n"}
